# Update RF values for rows 31-58 (Catch_Trust_31) from 218.2057 to 36.6954
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I31:I58").Value = 36.6954
